$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new extraction-summary row (row 17) below the existing data.
# Columns A-E hold the library name / raw counts (text + numbers); F-J
# hold percentage and timestamp values that must remain TEXT (like every
# other row already on the sheet) rather than being auto-coerced to
# numbers by Excel. Temporarily mark those cells as Text ("@") before
# writing the numeric-looking strings, then restore the default "Normal"
# style so the new cells don't carry a stray number-format style.
$ws.Range("A17").Value = "test"
$ws.Range("B17").Value = 2366
$ws.Range("C17").Value = 1229
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 1633

$ws.Range("F17:J17").NumberFormat = "@"
$ws.Range("F17").Value = "45.22"
$ws.Range("G17").Value = "23.49"
$ws.Range("H17").Value = "0.08"
$ws.Range("I17").Value = "31.21"
$ws.Range("J17").Value = "2025-08-29 18:36:58"
$ws.Range("F17:J17").Style = "Normal"
